# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.001754667048134761, 0.05231270169004087, 3.082599426703578, 0.4998867070740569, 0, 3.636553502515811)
    3 = @(1.505614041169197, 0.3375848360084654, 0.1529057820181812, 0.4998867070740569, 0, 2.495991366269901)
    4 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 6.48142807727062, 0, 12.0302756157461)
    5 = @(0.06328177979961902, 0.05231270169004087, 0.1529057820181812, 0.4998867070740569, 1, 0.768386970581898)
    6 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 0, 8.418600821238126)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("F$row").Value = $vals[4]
    $ws.Range("G$row").Value = $vals[5]
}
